$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.890.02'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.873.71'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7373'
$ws.Range('E5').Value = '  -4.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.30'
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07200'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('E10').Value = '  -4.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08406'
$ws.Range('E11').Value = '  -3.66%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7493'
$ws.Range('E12').Value = '  -3.13%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.888.55'
$ws.Range('E13').Value = '  -5.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.412'
$ws.Range('E14').Value = '  -0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.42'
$ws.Range('E15').Value = '  -2.19%  '
$ws.Range('D16').Value = '29.904.18'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.104'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.56'
$ws.Range('E18').Value = '  -2.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.06'
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007808'
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9993'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = '2.124.17'
$ws.Range('E22').Value = '  -8.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.980'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9996'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1548'
$ws.Range('E25').Value = '  -3.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.254'
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.64'
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.58'
$ws.Range('E28').Value = '  -1.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.028'
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.496'
$ws.Range('E30').Value = '  +4.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.574'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.532'
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.245'
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05312'
$ws.Range('E34').Value = '  -2.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.234'
$ws.Range('E35').Value = '  -1.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7525'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.695'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.752'
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4510'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('D42').Value = '1.113.97'
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.041'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.25'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8553'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.53'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.634'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.077'
$ws.Range('E49').Value = '  +2.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.835'
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('D51').Value = '2.021.75'
$ws.Range('E51').Value = '  -6.53%  '
